# it-180203 fix add sub-level with and without siblings
#
# Adds a new sub-level row (row 3) to the test_hierarchy_table sheet,
# mirroring row 2 (a second sibling under the same parent), resizes the
# table to include it, bumps row 2's height, and updates the active-sheet /
# selection bookkeeping (moving from monthly_issues -> test_hierarchy_table).

$wb = $excel.ActiveWorkbook

$wsHierarchy = $wb.Worksheets.Item("test_hierarchy_table")
$wsMonthly   = $wb.Worksheets.Item("monthly_issues")

# --- add the new row, seeded from row 2 (same parent/level) so it picks up
# the same number formats / alignment, then fix up the values that differ
# (seq goes to 2 for this second sibling under the same parent).
$wsHierarchy.Range("A2:G2").Copy($wsHierarchy.Range("A3:G3"))
$wsHierarchy.Range("B3").Value = 2

# --- grow the table so the new row is part of Table110 -------------------
$tbl = $wsHierarchy.ListObjects.Item(1)
$tbl.Resize($wsHierarchy.Range("A1:G3"))

# --- row 2 grows a touch taller now that there's a sibling below it ------
$wsHierarchy.Rows.Item(2).RowHeight = 33

# --- selection / active-sheet bookkeeping ---------------------------------
# monthly_issues was the active tab before; it loses that status and its
# selection moves on.
$wsMonthly.Activate()
$wsMonthly.Range("A12").Select()

# test_hierarchy_table becomes the active tab, selection on the new row's
# rgt cell.
$wsHierarchy.Activate()
$wsHierarchy.Range("D2").Select()
